# Wallkill_2018_flow_fieldsheets.xlsx — "Creating and updating plots"
#
# 1) Rename the "Flow" sheet to "flow" (lowercase). Excel automatically
#    updates the _xlnm._FilterDatabase defined name and the pivot cache's
#    worksheetSource reference that point at this sheet.
# 2) Update the header row (row 1) on that sheet: rename a few columns to
#    their lowercase/renamed forms.
# 3) Make "flow" the active sheet/tab, with cell E6 selected.
# 4) The "simplified qual" sheet stops being the active tab; its selection
#    moves to D19.

$wb = $excel.ActiveWorkbook

$flow = $wb.Worksheets.Item("Flow")
$flow.Name = "flow"

# Update header row text (order matches the original authoring sequence so
# newly-introduced shared strings land at the same indices as the target:
# "flow_quant_cfs" was introduced before "station").
$flow.Range("B1").Value = "sample_date"
$flow.Range("C1").Value = "Crew"
$flow.Range("D1").Value = "flow_qual"
$flow.Range("E1").Value = "flow_quant_cfs"
$flow.Range("F1").Value = "Uncertainty"
$flow.Range("G1").Value = "Lgst_error_src"
$flow.Range("H1").Value = "Highest_error_src_pct_err"
$flow.Range("I1").Value = "Notes"
$flow.Range("A1").Value = "station"

# Update the "simplified qual" sheet's selection before switching tabs
$qual = $wb.Worksheets.Item("simplified qual")
$qual.Range("D19").Select()

# Make "flow" active and select E6 (within the frozen header pane)
$flow.Activate()
$flow.Range("E6").Select()
